$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value corrections (error calculations) ---
$ws.Range("D3").Value = -14.2
$ws.Range("F4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("F9").Value = 17.26
$ws.Range("F10").Value = 16.43
$ws.Range("F13").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("D21").Value = -14.3
$ws.Range("D23").ClearContents()
$ws.Range("D34").Value = -14.7

# --- Remove rows that are no longer part of the data set ---
# "RM 232" (row 26) and "SC 92" (row 28) are dropped entirely; every row
# below shifts up by two. Delete bottom-up so row 26's index is still valid
# when it is removed second.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()
